$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FromLocation value on row 2 (column D) from "nyc" to "las"
$ws.Range("D2").Value = "las"

# Update the active selection to D13 to match the saved view state
$ws.Range("D13").Select()
